$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TEXTO")

# Update the label in B3 ("Error=" -> " RH="); every downstream formula
# (B5, row 8 CHAR, row 9 CODE, row 10 DEC2HEX, row 11/12 word values)
# recalculates automatically from this single change.
$ws.Range("B3").Value = " RH="

# Update the active selection on the TEXTO sheet to C5.
[void]$ws.Range("C5").Select()
